$wb = $excel.ActiveWorkbook

# --- Sheet "Orders": append new rows 62-70 in column C (PackageID strings) ---
$ordersWs = $wb.Worksheets.Item("Orders")

$newItems = @(
    "435_猫眼_Echinacea_undefined_1bunch",
    "401_大飞燕白色_delphinium white_undefined_1bunch",
    "403_大飞燕浅蓝色_delphinium light blue_undefined_1bunch",
    "496_大飞燕深蓝色_delphinium dark blue_undefined_1bunch",
    "571_大飞燕浅紫_undefined_undefined_1bunch",
    "402_大飞燕深紫色_delphinium purple_undefined_1bunch",
    "495_大飞燕深粉色_delphinium pink_undefined_1bunch",
    "433_红豆_Hypericum red_undefined_1bunch",
    "565_千日红_Gomphrena_undefined_1bunch"
)

$startRow = 62
for ($i = 0; $i -lt $newItems.Length; $i++) {
    $row = $startRow + $i
    $ordersWs.Cells.Item($row, 3).Value = $newItems[$i]
}

# Row 68 also carries an explicit (empty) value in column A per the source data.
# A bare "'" (quote-prefix) is the only way to make Excel persist a literal
# empty-string cell instead of clearing/removing it.
$ordersWs.Cells.Item(68, 1).Value = "'"

# --- Sheet "Summary": extend the digit-string in G2 with 9 more trailing zeros ---
$summaryWs = $wb.Worksheets.Item("Summary")
$g2 = $summaryWs.Cells.Item(2, 7)
$g2.NumberFormat = "@"
$g2.Value = "05200000000000000000000000000000000000000000000000000000000000000000000"
